$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0.9130434989929199
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0.5

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.5

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0.5

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0.5

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0.97826087474823
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0.489130437374115

$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0.97826087474823
$ws.Range("D9").Value = 0.97826087474823
$ws.Range("E9").Value = 0.9890105128288269
$ws.Range("F9").Value = 0.489130437374115

$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0.9999995231628418
$ws.Range("F10").Value = 0.5

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0.5

$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0.5

$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0.5

$ws.Range("B14").Value = 0.3999999165534973
$ws.Range("C14").Value = 0.2499999701976776
$ws.Range("D14").Value = 0.9021739363670349
$ws.Range("E14").Value = 0.3076917827129364
$ws.Range("F14").Value = 0.6071428656578064

$ws.Range("B15").Value = 0.9999997019767761
$ws.Range("C15").Value = 0.1874999850988388
$ws.Range("D15").Value = 0.8586956262588501
$ws.Range("E15").Value = 0.3157891929149628
$ws.Range("F15").Value = 0.59375

$ws.Range("B16").Value = 0.9999999403953552
$ws.Range("C16").Value = 0.7333332896232605
$ws.Range("D16").Value = 0.9130434989929199
$ws.Range("E16").Value = 0.8461533188819885
$ws.Range("F16").Value = 0.8666666746139526

$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0.9999995231628418
$ws.Range("F17").Value = 1

$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0.95652174949646
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0.47826087474823

$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0.5

$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0.989130437374115
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0.4945652186870575

$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0.5

$ws.Range("B22").Value = 0.9999998807907104
$ws.Range("C22").Value = 0.8749998807907104
$ws.Range("D22").Value = 0.989130437374115
$ws.Range("E22").Value = 0.9333327412605286
$ws.Range("F22").Value = 0.9374999403953552

$ws.Range("B23").Value = 0.08695652335882187
$ws.Range("C23").Value = 0.9999998807907104
$ws.Range("D23").Value = 0.08695652335882187
$ws.Range("E23").Value = 0.1599998623132706
$ws.Range("F23").Value = 0.4999999403953552

$ws.Range("D24").Value = 1

$ws.Range("D25").Value = 0.9130434989929199

$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0.08695652335882187
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 0.4999999403953552

$ws.Range("B27").Value = 0.4239130318164825
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 0.4239130318164825
$ws.Range("E27").Value = 0.5954194664955139
$ws.Range("F27").Value = 0.5

$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 0.4239130318164825
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0.5

$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0.5

$ws.Range("B30").Value = 0
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 0.9130434989929199
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0.5

$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 0.989130437374115
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0.4945652186870575

$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 0.95652174949646
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0.47826087474823

$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 0.72826087474823
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 0.364130437374115

$ws.Range("B34").Value = 0.04999999701976776
$ws.Range("C34").Value = 0.3333332240581512
$ws.Range("D34").Value = 0.77173912525177
$ws.Range("E34").Value = 0.08695628494024277
$ws.Range("F34").Value = 0.5599250197410583

$ws.Range("B35").Value = 0.07692307233810425
$ws.Range("C35").Value = 0.07692307233810425
$ws.Range("D35").Value = 0.47826087474823
$ws.Range("E35").Value = 0.07692257314920425
$ws.Range("F35").Value = 0.3566433489322662

$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 0.8913043737411499
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0.488095223903656

$ws.Range("D37").Value = 1

$ws.Range("B38").Value = 0.2857142686843872
$ws.Range("C38").Value = 0.08510638028383255
$ws.Range("D38").Value = 0.4239130318164825
$ws.Range("E38").Value = 0.1311471909284592
$ws.Range("F38").Value = 0.4314420819282532

$ws.Range("B39").Value = 0.8152173757553101
$ws.Range("C39").Value = 1
$ws.Range("D39").Value = 0.8152173757553101
$ws.Range("E39").Value = 0.8982031345367432
$ws.Range("F39").Value = 0.5

$ws.Range("B40").Value = 0
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 0.8260869383811951
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 0.5

$ws.Range("B41").Value = 0
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 0.989130437374115
$ws.Range("E41").Value = 0
$ws.Range("F41").Value = 0.5

$ws.Range("B42").Value = 0.555555522441864
$ws.Range("C42").Value = 0.8823528289794922
$ws.Range("D42").Value = 0.8478260636329651
$ws.Range("E42").Value = 0.6818176507949829
$ws.Range("F42").Value = 0.8611763715744019

$ws.Range("B43").Value = 0.7777776718139648
$ws.Range("C43").Value = 0.8749998807907104
$ws.Range("D43").Value = 0.967391312122345
$ws.Range("E43").Value = 0.8235287666320801
$ws.Range("F43").Value = 0.9255951642990112

$ws.Range("B44").Value = 0.8461537957191467
$ws.Range("C44").Value = 0.3793103098869324
$ws.Range("D44").Value = 0.782608687877655
$ws.Range("E44").Value = 0.5238090753555298
$ws.Range("F44").Value = 0.6737821102142334

$ws.Range("B45").Value = 0.8372092843055725
$ws.Range("C45").Value = 0.9473684430122375
$ws.Range("D45").Value = 0.9021739363670349
$ws.Range("E45").Value = 0.888888418674469
$ws.Range("F45").Value = 0.9088693857192993

$ws.Range("B46").Value = 0
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 1
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0.5

$ws.Range("B47").Value = 1
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = 1
$ws.Range("E47").Value = 0.9999995231628418
$ws.Range("F47").Value = 0.5

$ws.Range("B48").Value = 0
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 1
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0.5

$ws.Range("D49").Value = 1
$ws.Range("F49").Value = 0.5
